$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / unambiguous updates (coin names, links, prices with thousands separators, percentages) ---
$ws.Range("D2").Value = "61.664.57"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "3.002.02"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "2.999.96"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "3.498.38"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "61.660.61"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "3.006.43"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E25").Value = "  +10.85%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +3.86%  "
$ws.Range("D35").Value = "0.0₃0836"
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +11.16%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "2.706.50"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +2.69%  "

# --- Numeric-looking price updates: force text so Excel does not coerce them to numbers ---
$numericLooking = @("D5","D6","D12","D13","D14","D17","D20","D21","D23","D24","D25","D26","D27","D30","D33","D40","D41","D42","D43","D44","D45","D48","D51")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "599.95"
$ws.Range("D6").Value = "145.15"
$ws.Range("D12").Value = "0.459"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D14").Value = "34.54"
$ws.Range("D17").Value = "7.01"
$ws.Range("D20").Value = "452.19"
$ws.Range("D21").Value = "14.04"
$ws.Range("D23").Value = "7.36"
$ws.Range("D24").Value = "81.72"
$ws.Range("D25").Value = "10.98"
$ws.Range("D26").Value = "2.27"
$ws.Range("D27").Value = "12.04"
$ws.Range("D30").Value = "7.28"
$ws.Range("D33").Value = "27.50"
$ws.Range("D40").Value = "50.47"
$ws.Range("D41").Value = "0.125"
$ws.Range("D42").Value = "2.91"
$ws.Range("D43").Value = "403.92"
$ws.Range("D44").Value = "39.72"
$ws.Range("D45").Value = "0.274"
$ws.Range("D48").Value = "131.88"
$ws.Range("D51").Value = "2.16"
foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
